$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (89 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 990.7222
$ws.Range("J17").Value = 990.7222
$ws.Range("L17").Value = 2972.1666
$ws.Range("N17").Value = -3308.1666
$ws.Range("H64").Value = 7386.6333
$ws.Range("J64").Value = 8271.959999999999
$ws.Range("L64").Value = 8271.959999999999
$ws.Range("N64").Value = -8767.959999999999
$ws.Range("H67").Value = 7386.6333
$ws.Range("J67").Value = 8271.959999999999
$ws.Range("L67").Value = 8271.959999999999
$ws.Range("N67").Value = -9987.959999999999
$ws.Range("H70").Value = 46160940
$ws.Range("I70").Value = 20004766
$ws.Range("J70").Value = 62508550
$ws.Range("K70").Value = 60014298
$ws.Range("L70").Value = 187525650
$ws.Range("M70").Value = -60014028
$ws.Range("N70").Value = -187526190
$ws.Range("H73").Value = 46160940
$ws.Range("I73").Value = 20004766
$ws.Range("J73").Value = 62508550
$ws.Range("K73").Value = 60014298
$ws.Range("L73").Value = 187525650
$ws.Range("M73").Value = -60013362
$ws.Range("N73").Value = -187527522
$ws.Range("H76").Value = 2543.4167
$ws.Range("J76").Value = 2736.5
$ws.Range("L76").Value = 2736.5
$ws.Range("N76").Value = -3366.5
$ws.Range("H79").Value = 2543.4167
$ws.Range("J79").Value = 2736.5
$ws.Range("L79").Value = 2736.5
$ws.Range("N79").Value = -4920.5
$ws.Range("H97").Value = 5098.3
$ws.Range("J97").Value = 5098.3
$ws.Range("L97").Value = 15294.9
$ws.Range("N97").Value = -16286.9
$ws.Range("H99").Value = 1323.5385
$ws.Range("I99").Value = 655.1818
$ws.Range("J99").Value = 4999.5
$ws.Range("K99").Value = 1965.5454
$ws.Range("L99").Value = 14998.5
$ws.Range("M99").Value = -467.5454
$ws.Range("N99").Value = -17994.5
$ws.Range("H107").Value = 8877.111000000001
$ws.Range("I107").Value = 9818.093999999999
$ws.Range("K107").Value = 9818.093999999999
$ws.Range("M107").Value = -7898.093999999999
$ws.Range("H115").Value = 402.5
$ws.Range("I115").Value = 286.66666
$ws.Range("K115").Value = 859.9999799999999
$ws.Range("M115").Value = 707.0000200000001
$ws.Range("H118").Value = 1738.2667
$ws.Range("J118").Value = 2920
$ws.Range("L118").Value = 8760
$ws.Range("N118").Value = -12074
$ws.Range("H125").Value = 14558.588
$ws.Range("I125").Value = 3032
$ws.Range("J125").Value = 27526
$ws.Range("K125").Value = 27288
$ws.Range("L125").Value = 247734
$ws.Range("M125").Value = -24828
$ws.Range("N125").Value = -252654
$ws.Range("H131").Value = 6896.5557
$ws.Range("I131").Value = 2937.5
$ws.Range("K131").Value = 8812.5
$ws.Range("M131").Value = -3772.5
$ws.Range("H135").Value = 1143.6316
$ws.Range("I135").Value = 1043.7333
$ws.Range("J135").Value = 1518.25
$ws.Range("K135").Value = 9393.599700000001
$ws.Range("L135").Value = 13664.25
$ws.Range("M135").Value = -6858.599700000001
$ws.Range("N135").Value = -18734.25
$ws.Range("H138").Value = 4340
$ws.Range("I138").Value = 3707.2307
$ws.Range("J138").Value = 4797
$ws.Range("K138").Value = 11121.6921
$ws.Range("L138").Value = 14391
$ws.Range("M138").Value = -5981.6921
$ws.Range("N138").Value = -24671
$ws.Range("H141").Value = 2811.5293
$ws.Range("I141").Value = 2292
$ws.Range("J141").Value = 4500
$ws.Range("K141").Value = 6876
$ws.Range("L141").Value = 13500
$ws.Range("M141").Value = -1696
$ws.Range("N141").Value = -23860

# --- Sheet: ARM (49 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3125.5625
$ws.Range("J32").Value = 3498.6667
$ws.Range("L32").Value = 3498.6667
$ws.Range("N32").Value = -4072.6667
$ws.Range("H61").Value = 5284.826
$ws.Range("I61").Value = 5238.8423
$ws.Range("K61").Value = 5238.8423
$ws.Range("M61").Value = -5026.8423
$ws.Range("H74").Value = 1559.275
$ws.Range("I74").Value = 1615.9429
$ws.Range("J74").Value = 1162.6
$ws.Range("K74").Value = 1615.9429
$ws.Range("L74").Value = 1162.6
$ws.Range("M74").Value = -741.9429
$ws.Range("N74").Value = -2910.6
$ws.Range("H77").Value = 1559.275
$ws.Range("I77").Value = 1615.9429
$ws.Range("J77").Value = 1162.6
$ws.Range("K77").Value = 8079.7145
$ws.Range("L77").Value = 5813
$ws.Range("M77").Value = -3711.7145
$ws.Range("N77").Value = -14549
$ws.Range("H109").Value = 105188.5
$ws.Range("J109").Value = 105188.5
$ws.Range("L109").Value = 105188.5
$ws.Range("N109").Value = -107962.5
$ws.Range("H110").Value = 921.75
$ws.Range("I110").Value = 921.75
$ws.Range("K110").Value = 921.75
$ws.Range("M110").Value = 1123.25
$ws.Range("H122").Value = 2927254
$ws.Range("I122").Value = 4276813
$ws.Range("J122").Value = 3209.0833
$ws.Range("K122").Value = 12830439
$ws.Range("L122").Value = 9627.249899999999
$ws.Range("M122").Value = -12827989
$ws.Range("N122").Value = -14527.2499
$ws.Range("H132").Value = 2753.2222
$ws.Range("I132").Value = 2621.0881
$ws.Range("K132").Value = 7863.2643
$ws.Range("M132").Value = -5333.2643
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140
$ws.Range("H136").Value = 5284.826
$ws.Range("I136").Value = 5238.8423
$ws.Range("K136").Value = 15716.5269
$ws.Range("M136").Value = -13166.5269

# --- Sheet: BSM (19 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 363.80768
$ws.Range("J80").Value = 374.91666
$ws.Range("L80").Value = 374.91666
$ws.Range("N80").Value = -2370.91666
$ws.Range("H83").Value = 363.80768
$ws.Range("J83").Value = 374.91666
$ws.Range("L83").Value = 1874.5833
$ws.Range("N83").Value = -11858.5833
$ws.Range("H123").Value = 117499.5
$ws.Range("J123").Value = 117499.5
$ws.Range("L123").Value = 117499.5
$ws.Range("N123").Value = -127299.5
$ws.Range("H134").Value = 3801.1904
$ws.Range("I134").Value = 4042.0667
$ws.Range("J134").Value = 3199
$ws.Range("K134").Value = 12126.2001
$ws.Range("L134").Value = 9597
$ws.Range("M134").Value = -9591.2001
$ws.Range("N134").Value = -14667

# --- Sheet: CRP (43 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2599
$ws.Range("I2").Value = 199.5
$ws.Range("J2").Value = 4998.5
$ws.Range("K2").Value = 199.5
$ws.Range("L2").Value = 4998.5
$ws.Range("M2").Value = -86.5
$ws.Range("N2").Value = -5224.5
$ws.Range("H31").Value = 4254.7085
$ws.Range("I31").Value = 1692.5294
$ws.Range("J31").Value = 10477.143
$ws.Range("K31").Value = 1692.5294
$ws.Range("L31").Value = 10477.143
$ws.Range("M31").Value = -1397.5294
$ws.Range("N31").Value = -11067.143
$ws.Range("H34").Value = 4254.7085
$ws.Range("I34").Value = 1692.5294
$ws.Range("J34").Value = 10477.143
$ws.Range("K34").Value = 1692.5294
$ws.Range("L34").Value = 10477.143
$ws.Range("M34").Value = -1490.5294
$ws.Range("N34").Value = -10881.143
$ws.Range("H94").Value = 2372.3076
$ws.Range("I94").Value = 2126.6667
$ws.Range("J94").Value = 2582.8572
$ws.Range("K94").Value = 2126.6667
$ws.Range("L94").Value = 2582.8572
$ws.Range("M94").Value = -1675.6667
$ws.Range("N94").Value = -3484.8572
$ws.Range("H122").Value = 2524.5898
$ws.Range("I122").Value = 2385.077
$ws.Range("K122").Value = 7155.231000000001
$ws.Range("M122").Value = -4705.231000000001
$ws.Range("H132").Value = 2236.2273
$ws.Range("I132").Value = 1789.75
$ws.Range("J132").Value = 6701
$ws.Range("K132").Value = 5369.25
$ws.Range("L132").Value = 20103
$ws.Range("M132").Value = -2839.25
$ws.Range("N132").Value = -25163
$ws.Range("H134").Value = 1857.3334
$ws.Range("I134").Value = 1848.9412
$ws.Range("K134").Value = 5546.8236
$ws.Range("M134").Value = -3011.8236

# --- Sheet: CUL (51 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2184.3845
$ws.Range("I3").Value = 2529.8
$ws.Range("K3").Value = 7589.400000000001
$ws.Range("M3").Value = -7477.400000000001
$ws.Range("H6").Value = 78.875
$ws.Range("I6").Value = 96.5
$ws.Range("K6").Value = 289.5
$ws.Range("M6").Value = -176.5
$ws.Range("H18").Value = 2052
$ws.Range("I18").Value = 390.53845
$ws.Range("J18").Value = 4451.8887
$ws.Range("K18").Value = 1171.61535
$ws.Range("L18").Value = 13355.6661
$ws.Range("M18").Value = -1002.61535
$ws.Range("N18").Value = -13693.6661
$ws.Range("H46").Value = 9989112
$ws.Range("I46").Value = 23299288
$ws.Range("J46").Value = 3334024.8
$ws.Range("K46").Value = 69897864
$ws.Range("L46").Value = 10002074.4
$ws.Range("M46").Value = -69897773
$ws.Range("N46").Value = -10002256.4
$ws.Range("H68").Value = 3454.625
$ws.Range("J68").Value = 3845.6365
$ws.Range("L68").Value = 11536.9095
$ws.Range("N68").Value = -13158.9095
$ws.Range("H71").Value = 3454.625
$ws.Range("J71").Value = 3845.6365
$ws.Range("L71").Value = 34610.7285
$ws.Range("N71").Value = -42722.7285
$ws.Range("H134").Value = 3935.7083
$ws.Range("I134").Value = 1859.7273
$ws.Range("J134").Value = 5692.3076
$ws.Range("K134").Value = 5579.1819
$ws.Range("L134").Value = 17076.9228
$ws.Range("M134").Value = -509.1818999999996
$ws.Range("N134").Value = -27216.9228
$ws.Range("H139").Value = 2715.3333
$ws.Range("I139").Value = 2231.5557
$ws.Range("J139").Value = 4166.6665
$ws.Range("K139").Value = 6694.6671
$ws.Range("L139").Value = 12499.9995
$ws.Range("M139").Value = -1554.6671
$ws.Range("N139").Value = -22779.9995
$ws.Range("H140").Value = 3480.85
$ws.Range("I140").Value = 2601.1875
$ws.Range("J140").Value = 6999.5
$ws.Range("K140").Value = 7803.5625
$ws.Range("L140").Value = 20998.5
$ws.Range("M140").Value = -2623.5625
$ws.Range("N140").Value = -31358.5

# --- Sheet: GSM (41 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11121845
$ws.Range("J70").Value = 11103.345
$ws.Range("L70").Value = 11103.345
$ws.Range("N70").Value = -11643.345
$ws.Range("H73").Value = 11121845
$ws.Range("J73").Value = 11103.345
$ws.Range("L73").Value = 11103.345
$ws.Range("N73").Value = -12975.345
$ws.Range("H80").Value = 12481.23
$ws.Range("I80").Value = 18775
$ws.Range("J80").Value = 7086.5713
$ws.Range("K80").Value = 18775
$ws.Range("L80").Value = 7086.5713
$ws.Range("M80").Value = -17777
$ws.Range("N80").Value = -9082.5713
$ws.Range("H83").Value = 12481.23
$ws.Range("I83").Value = 18775
$ws.Range("J83").Value = 7086.5713
$ws.Range("K83").Value = 93875
$ws.Range("L83").Value = 35432.85649999999
$ws.Range("M83").Value = -88883
$ws.Range("N83").Value = -45416.85649999999
$ws.Range("H102").Value = 2086.6365
$ws.Range("I102").Value = 2086.6365
$ws.Range("K102").Value = 2086.6365
$ws.Range("M102").Value = -464.6365000000001
$ws.Range("H117").Value = 75000
$ws.Range("J117").Value = 75000
$ws.Range("L117").Value = 75000
$ws.Range("N117").Value = -81884
$ws.Range("H122").Value = 9777.162
$ws.Range("I122").Value = 10786.033
$ws.Range("K122").Value = 32358.099
$ws.Range("M122").Value = -29908.099
$ws.Range("H132").Value = 3396.5454
$ws.Range("I132").Value = 3576.6875
$ws.Range("J132").Value = 2916.1667
$ws.Range("K132").Value = 10730.0625
$ws.Range("L132").Value = 8748.500100000001
$ws.Range("M132").Value = -8200.0625
$ws.Range("N132").Value = -13808.5001

# --- Sheet: LTW (19 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2706.5833
$ws.Range("I46").Value = 2123.75
$ws.Range("K46").Value = 2123.75
$ws.Range("M46").Value = -1935.75
$ws.Range("H122").Value = 5474.727
$ws.Range("I122").Value = 5049.8125
$ws.Range("J122").Value = 6607.8335
$ws.Range("K122").Value = 15149.4375
$ws.Range("L122").Value = 19823.5005
$ws.Range("M122").Value = -12699.4375
$ws.Range("N122").Value = -24723.5005
$ws.Range("H132").Value = 4112.478
$ws.Range("I132").Value = 4566.067
$ws.Range("K132").Value = 13698.201
$ws.Range("M132").Value = -11168.201
$ws.Range("H136").Value = 6099.4165
$ws.Range("I136").Value = 5319.6
$ws.Range("K136").Value = 15958.8
$ws.Range("M136").Value = -13408.8

# --- Sheet: WVR (48 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 9990
$ws.Range("J45").Value = 9980
$ws.Range("L45").Value = 9980
$ws.Range("N45").Value = -10962
$ws.Range("H74").Value = 31311.5
$ws.Range("J74").Value = 28082
$ws.Range("L74").Value = 28082
$ws.Range("N74").Value = -29954
$ws.Range("H77").Value = 31311.5
$ws.Range("J77").Value = 28082
$ws.Range("L77").Value = 84246
$ws.Range("N77").Value = -93606
$ws.Range("H100").Value = 3910
$ws.Range("J100").Value = 4000.5
$ws.Range("L100").Value = 8001
$ws.Range("N100").Value = -9083
$ws.Range("H107").Value = 2401.6538
$ws.Range("I107").Value = 2503.0476
$ws.Range("J107").Value = 1975.8
$ws.Range("K107").Value = 7509.1428
$ws.Range("L107").Value = 5927.4
$ws.Range("M107").Value = -5589.1428
$ws.Range("N107").Value = -9767.4
$ws.Range("H122").Value = 6603
$ws.Range("I122").Value = 6085
$ws.Range("J122").Value = 7121
$ws.Range("K122").Value = 18255
$ws.Range("L122").Value = 21363
$ws.Range("M122").Value = -15805
$ws.Range("N122").Value = -26263
$ws.Range("H132").Value = 4310.3687
$ws.Range("I132").Value = 3545.2222
$ws.Range("J132").Value = 6188.4546
$ws.Range("K132").Value = 10635.6666
$ws.Range("L132").Value = 18565.3638
$ws.Range("M132").Value = -8105.6666
$ws.Range("N132").Value = -23625.3638
$ws.Range("H136").Value = 1802.5927
$ws.Range("I136").Value = 1723.4762
$ws.Range("J136").Value = 2079.5
$ws.Range("K136").Value = 5170.4286
$ws.Range("L136").Value = 6238.5
$ws.Range("M136").Value = -2620.4286
$ws.Range("N136").Value = -11338.5
$ws.Range("H138").Value = 52486.5
$ws.Range("J138").Value = 52486.5
$ws.Range("L138").Value = 52486.5
$ws.Range("N138").Value = -62766.5

Write-Host "Applied all cell updates."